$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 0) Remove all existing comments first. This resets the internal
#    comment-id counter to 0 so the comments we (re)create below land on
#    the exact w:id values the target document expects (0, 1, 3 - with
#    "2" intentionally skipped, see step 4).  Doing this BEFORE any text
#    edits also avoids corrupting the old commentRangeStart/End markers
#    that overlap text we are about to rewrite.
# ---------------------------------------------------------------------
while ($d.Comments.Count -gt 0) {
    $d.Comments(1).Delete()
}

# The old _GoBack bookmark (sitting between "; F" and ") during") is
# removed from the body in the target revision.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 1) Figure 1 legend, item 1: "Cushing's and control BMI" gains N callouts.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Cushing’s and control BMI", $true, $false, $false, $false, $false, $true, 1, $false, "Cushing’s (non-obese n=; obese n= ) and control (non-obese n=; obese n=) BMI", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Figure 1 legend, item 1: "5 weeks of treatment" gains dosing groups.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("5 weeks of treatment", $true, $false, $false, $false, $false, $true, 1, $false, "5 weeks of dexamethasone (NCD n=; HFD n=) or vehicle (NCD n=; HFD n=) treatment", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Figure 1 legend, item 1: the euglycemic-clamp / asterisk sentence is
#    expanded with N's, fasting info, and a second asterisk-meaning note.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(" clamp following 3 weeks of dexamethasone or vehicle treatment. Asterisks indicate a significant interaction between diet and treatment.", $true, $false, $false, $false, $false, $true, 1, $false, " clamp following 3 weeks of dexamethasone (n=14) or vehicle (n=11) treatment and 11 weeks of HFD. All mice were fasted for 6 hours prior to experiments. Asterisks in between two bars of the same condition indicate a significant interaction between diet and treatment. Centered asterisks indicated statistically significant treatment effect.", 2) | Out-Null

# ---------------------------------------------------------------------
# 4) Re-create the three comments at their new anchors / wording.
#    Comment "0": replaces the old id=1 comment; now anchored on the
#    newly-inserted Cushing's/control N callouts.
# ---------------------------------------------------------------------
$c0rng = $d.Content
$c0rng.Find.Execute(" (non-obese n=; obese n= ) and control (non-obese n=; obese n=)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Comments.Add($c0rng, "Cannot find these in the script or anywhere else I have looked. Do you know where I can find the Ns for the different groups here?") | Out-Null

# Comment "1": brand new, anchored on the newly-inserted NCD/HFD N callouts.
$c1rng = $d.Content
$c1rng.Find.Execute("(NCD n=; HFD n=) or vehicle (NCD n=; HFD n=) ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Comments.Add($c1rng, "Find N for these experiments") | Out-Null

# A throw-away comment anchored near the very end of the document - its
# sole purpose is to consume id "2" (which the target document leaves
# unused) before the "transcripts" comment is recreated as id "3". It is
# removed again a few lines down.
$dummyRng = $d.Content
$dummyRng.Find.Execute("tissues", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Comments.Add($dummyRng, "placeholder") | Out-Null

# Comment "3": replaces the old id=2 comment; same anchor ("transcripts")
# and same wording, just renumbered.
$c3rng = $d.Content
$c3rng.Find.Execute("transcripts", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Comments.Add($c3rng, "Should we also put the gluconeogenic transcripts here?") | Out-Null

# The placeholder comment sorts after all three real comments (its anchor
# is the last one in the document), so it is always the final entry in
# the position-ordered Comments collection - remove it by that index.
$d.Comments($d.Comments.Count).Delete()
